# Auto-generated script applying the Coeurl_Profits market-data refresh diff.
# For each (sheet, cell) pair: update the cached numeric value, add newly
# introduced cells, and clear cells that the diff removes entirely (clearing
# drops the backing <c> element on save, matching the source diff).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 13488.375
$ws.Range("I43").Value = 20979.6
$ws.Range("J43").Value = 1003
$ws.Range("K43").Value = 20979.6
$ws.Range("L43").Value = 1003
$ws.Range("M43").Value = -20910.6
$ws.Range("N43").Value = -1141
$ws.Range("H62").Value = 18186234
$ws.Range("I62").Value = 20004458
$ws.Range("K62").Value = 20004458
$ws.Range("M62").Value = -20003834
$ws.Range("H64").Value = 10849.2
$ws.Range("J64").Value = 13899.714
$ws.Range("L64").Value = 13899.714
$ws.Range("N64").Value = -14395.714
$ws.Range("H65").Value = 18186234
$ws.Range("I65").Value = 20004458
$ws.Range("K65").Value = 100022290
$ws.Range("M65").Value = -100019170
$ws.Range("H67").Value = 10849.2
$ws.Range("J67").Value = 13899.714
$ws.Range("L67").Value = 13899.714
$ws.Range("N67").Value = -15615.714
$ws.Range("H86").Value = 4326.909
$ws.Range("I86").Value = 2414.5715
$ws.Range("J86").Value = 7673.5
$ws.Range("K86").Value = 2414.5715
$ws.Range("L86").Value = 7673.5
$ws.Range("M86").Value = -1291.5715
$ws.Range("N86").Value = -9919.5
$ws.Range("H89").Value = 4326.909
$ws.Range("I89").Value = 2414.5715
$ws.Range("J89").Value = 7673.5
$ws.Range("K89").Value = 12072.8575
$ws.Range("L89").Value = 38367.5
$ws.Range("M89").Value = -6456.8575
$ws.Range("N89").Value = -49599.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4663.6445
$ws.Range("I32").Value = 4287.5
$ws.Range("K32").Value = 4287.5
$ws.Range("M32").Value = -4000.5
$ws.Range("H45").Value = 8673.895
$ws.Range("I45").Value = 8768.611000000001
$ws.Range("K45").Value = 8768.611000000001
$ws.Range("M45").Value = -8391.611000000001
$ws.Range("H92").Value = 29275
$ws.Range("I92").Value = 29000
$ws.Range("J92").Value = 29550
$ws.Range("K92").Value = 29000
$ws.Range("L92").Value = 29550
$ws.Range("M92").Value = -26504
$ws.Range("N92").Value = -34542

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1845.5294
$ws.Range("I86").Value = 1845.5294
$ws.Range("K86").Value = 1845.5294
$ws.Range("M86").Value = -722.5293999999999
$ws.Range("H89").Value = 1845.5294
$ws.Range("I89").Value = 1845.5294
$ws.Range("K89").Value = 9227.646999999999
$ws.Range("M89").Value = -3611.646999999999
$ws.Range("H132").Value = 59685
$ws.Range("J132").Value = 59685
$ws.Range("L132").Value = 59685
$ws.Range("N132").Value = -69805
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2454.238
$ws.Range("I31").Value = 1692.3
$ws.Range("K31").Value = 1692.3
$ws.Range("M31").Value = -1397.3
$ws.Range("H34").Value = 2454.238
$ws.Range("I34").Value = 1692.3
$ws.Range("K34").Value = 1692.3
$ws.Range("M34").Value = -1490.3
$ws.Range("H44").Value = 15000
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = ""
$ws.Range("H58").Value = 2654.5107
$ws.Range("I58").Value = 2807.8333
$ws.Range("J58").Value = 2383.9412
$ws.Range("K58").Value = 2807.8333
$ws.Range("L58").Value = 2383.9412
$ws.Range("M58").Value = -2604.8333
$ws.Range("N58").Value = -2789.9412
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = ""
$ws.Range("H62").Value = 3637.25
$ws.Range("I62").Value = 3949.5
$ws.Range("J62").Value = 3325
$ws.Range("K62").Value = 3949.5
$ws.Range("L62").Value = 3325
$ws.Range("M62").Value = -3325.5
$ws.Range("N62").Value = -4573
$ws.Range("H65").Value = 3637.25
$ws.Range("I65").Value = 3949.5
$ws.Range("J65").Value = 3325
$ws.Range("K65").Value = 19747.5
$ws.Range("L65").Value = 16625
$ws.Range("M65").Value = -16627.5
$ws.Range("N65").Value = -22865
$ws.Range("H107").Value = 743.89655
$ws.Range("J107").Value = 1419.1818
$ws.Range("L107").Value = 1419.1818
$ws.Range("N107").Value = -5259.1818
$ws.Range("H132").Value = 4129.778
$ws.Range("I132").Value = 3944.276
$ws.Range("K132").Value = 11832.828
$ws.Range("M132").Value = -9302.828
$ws.Range("H136").Value = 2654.5107
$ws.Range("I136").Value = 2807.8333
$ws.Range("J136").Value = 2383.9412
$ws.Range("K136").Value = 8423.499899999999
$ws.Range("L136").Value = 7151.823600000001
$ws.Range("M136").Value = -5873.499899999999
$ws.Range("N136").Value = -12251.8236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 258.36365
$ws.Range("J12").Value = 454.16666
$ws.Range("L12").Value = 1362.49998
$ws.Range("N12").Value = -1708.49998
$ws.Range("H25").Value = 110.5
$ws.Range("I25").Value = 110.5
$ws.Range("K25").Value = 331.5
$ws.Range("M25").Value = -162.5
$ws.Range("H30").Value = 110.5
$ws.Range("I30").Value = 110.5
$ws.Range("K30").Value = 331.5
$ws.Range("M30").Value = -229.5
$ws.Range("H40").Value = 329.25
$ws.Range("I40").Value = 177.28572
$ws.Range("J40").Value = 542
$ws.Range("K40").Value = 709.14288
$ws.Range("L40").Value = 2168
$ws.Range("M40").Value = -640.14288
$ws.Range("N40").Value = -2306
$ws.Range("H62").Value = 2019.5
$ws.Range("I62").Value = 723.4
$ws.Range("J62").Value = 8500
$ws.Range("K62").Value = 2170.2
$ws.Range("L62").Value = 25500
$ws.Range("M62").Value = -1484.2
$ws.Range("N62").Value = -26872
$ws.Range("H63").Value = 5000
$ws.Range("I63").Value = 5000
$ws.Range("K63").Value = 15000
$ws.Range("M63").Value = -14251
$ws.Range("H65").Value = 2019.5
$ws.Range("I65").Value = 723.4
$ws.Range("J65").Value = 8500
$ws.Range("K65").Value = 6510.599999999999
$ws.Range("L65").Value = 76500
$ws.Range("M65").Value = -3078.599999999999
$ws.Range("N65").Value = -83364
$ws.Range("H66").Value = 5000
$ws.Range("I66").Value = 5000
$ws.Range("K66").Value = 45000
$ws.Range("M66").Value = -41256
$ws.Range("H75").Value = 300
$ws.Range("I75").Value = 300
$ws.Range("K75").Value = 900
$ws.Range("M75").Value = 98
$ws.Range("H78").Value = 300
$ws.Range("I78").Value = 300
$ws.Range("K78").Value = 2700
$ws.Range("M78").Value = 2292
$ws.Range("H113").Value = 549.7619
$ws.Range("J113").Value = 619.26666
$ws.Range("L113").Value = 1857.79998
$ws.Range("N113").Value = -6197.79998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""
$ws.Range("H102").Value = 35715730
$ws.Range("I102").Value = 930.381
$ws.Range("J102").Value = 142860110
$ws.Range("K102").Value = 930.381
$ws.Range("L102").Value = 142860110
$ws.Range("M102").Value = 691.619
$ws.Range("N102").Value = -142863354
$ws.Range("H132").Value = 4116.4
$ws.Range("I132").Value = 3983.6365
$ws.Range("J132").Value = 4481.5
$ws.Range("K132").Value = 11950.9095
$ws.Range("L132").Value = 13444.5
$ws.Range("M132").Value = -9420.9095
$ws.Range("N132").Value = -18504.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9514.632
$ws.Range("I7").Value = 10906.692
$ws.Range("K7").Value = 10906.692
$ws.Range("M7").Value = -10794.692
$ws.Range("H68").Value = 3389.8
$ws.Range("I68").Value = 3000
$ws.Range("K68").Value = 3000
$ws.Range("M68").Value = -2251
$ws.Range("H71").Value = 3389.8
$ws.Range("I71").Value = 3000
$ws.Range("K71").Value = 15000
$ws.Range("M71").Value = -11256
$ws.Range("H126").Value = 9514.632
$ws.Range("I126").Value = 10906.692
$ws.Range("K126").Value = 32720.076
$ws.Range("M126").Value = -30250.076

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 169444.11
$ws.Range("I2").Value = 221666.33
$ws.Range("J2").Value = 64999.668
$ws.Range("K2").Value = 221666.33
$ws.Range("L2").Value = 64999.668
$ws.Range("M2").Value = -221554.33
$ws.Range("N2").Value = -65223.668
$ws.Range("H3").Value = 203279.8
$ws.Range("I3").Value = 1000000
$ws.Range("K3").Value = 1000000
$ws.Range("M3").Value = -999886
$ws.Range("H4").Value = 56312.22
$ws.Range("I4").Value = 126020.1
$ws.Range("J4").Value = 2690.7693
$ws.Range("K4").Value = 126020.1
$ws.Range("L4").Value = 2690.7693
$ws.Range("M4").Value = -125907.1
$ws.Range("N4").Value = -2916.7693
$ws.Range("H132").Value = 2096.742
$ws.Range("I132").Value = 2066.6333
$ws.Range("K132").Value = 6199.8999
$ws.Range("M132").Value = -3669.8999
